$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add the two sentences ("Kich thuoc cua bo nho dem ... theo hai cach:")
#    to the paragraph that currently only holds the _GoBack bookmark, placing
#    the new text *before* the bookmark so the bookmark stays at the end of
#    the paragraph (matching the diff).
# ---------------------------------------------------------------------------

$paras = $d.Paragraphs
$lastIndex = $paras.Count
$lastPara = $paras.Item($lastIndex)

# A run a couple of paragraphs earlier already carries the <w:lang val="en-US"/>
# run-formatting we need to reproduce; grab one character from it to use as a
# formatting template.
$templateRange = $paras.Item($lastIndex - 1).Range.Duplicate
$templateRange.SetRange($templateRange.Start, $templateRange.Start + 1)

$newSentence = "Kích thước của bộ nhớ đệm ảnh hưởng đến hiệu suất của External Sort theo hai cách:"

$bm = $d.Bookmarks.Item("_GoBack")
$insertPos = $bm.Range.Start

# Insert a one-character placeholder *before* the bookmark so the bookmark
# keeps sitting after our new content (InsertBefore on the bookmark range
# keeps the correct left-to-right order).
$bm.Range.InsertBefore("X")

# Stamp that placeholder character with the correct run formatting (lang=en-US)
# by copying the formatted text from the template range.
$anchor = $d.Range($insertPos, $insertPos + 1)
$anchor.FormattedText = $templateRange.FormattedText

# Now type the real sentence right after the freshly formatted placeholder -
# it inherits the placeholder run's formatting and merges into one run.
$anchor.Collapse(0)
$anchor.InsertAfter($newSentence)

# Remove the one-character placeholder.
$placeholder = $d.Range($insertPos, $insertPos + 1)
$placeholder.Delete()

# ---------------------------------------------------------------------------
# 2) Append two brand-new paragraphs after the (former) bookmark paragraph,
#    each holding one explanatory sentence.
# ---------------------------------------------------------------------------

$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$endRange = $lastPara.Range.Duplicate
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$paras = $d.Paragraphs
$secondPara = $paras.Item($paras.Count)
$secondRange = $secondPara.Range.Duplicate
$secondRange.Collapse(1)
$secondRange.InsertBefore("Tăng tốc độ sắp xếp: Khi kích thước của bộ nhớ đệm tăng, thuật toán sắp xếp sẽ có thể đọc và ghi dữ liệu nhiều hơn trong mỗi lần truy cập đĩa. Điều này sẽ giúp tăng tốc độ sắp xếp.")

$paras = $d.Paragraphs
$secondPara = $paras.Item($paras.Count)
$secondRange = $secondPara.Range.Duplicate
$secondRange.Collapse(0)
$secondRange.InsertParagraphAfter()

$paras = $d.Paragraphs
$thirdPara = $paras.Item($paras.Count)
$thirdRange = $thirdPara.Range.Duplicate
$thirdRange.Collapse(1)
$thirdRange.InsertBefore("Giảm số lượng truy cập đĩa: Khi kích thước của bộ nhớ đệm tăng, thuật toán sắp xếp sẽ cần ít truy cập đĩa hơn để sắp xếp dữ liệu. Điều này sẽ giúp giảm thời gian thực thi của thuật toán sắp xếp.")
